$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "content" column currently lives in column D (1-based col 4).
# We need to introduce a new "status" column at D, pushing the existing
# "content" column to E - without disturbing the column-width formatting
# that stays attached to column D.

# 1) Copy the existing column D ("content" header + values) over to column E.
$ws.Cells.Item(1, 5).Value = $ws.Cells.Item(1, 4).Value()
$ws.Cells.Item(2, 5).Value = $ws.Cells.Item(2, 4).Value()
$ws.Cells.Item(3, 5).Value = $ws.Cells.Item(3, 4).Value()

# 2) Overwrite column D with the new "status" data.
$ws.Cells.Item(1, 4).Value = "status"
$ws.Cells.Item(2, 4).Value = "Active"
$ws.Cells.Item(3, 4).Value = "Active"

# 3) Update the saved selection/active cell to match the new state.
$ws.Range("D6").Select() | Out-Null
